# THz-Compressed-Sensing metadata.xlsx update:
# add a new "1D NEW" table row (1dmockanderrors22.csv) describing a simulated
# interferogram run that resolves a second experimental peak.

$wb = $excel.ActiveWorkbook
$wsNew = $wb.Worksheets.Item("1D NEW")
$wsOld = $wb.Worksheets.Item("1D OLD")

# Grow the table ("Table43") by one row and fill it in.
$lo = $wsNew.ListObjects.Item(1)
$newRow = $lo.ListRows.Add()

$wsNew.Range("B25").Value = "1dmockanderrors22.csv"
$wsNew.Range("C25").Value = 53
$wsNew.Range("D25").Value = 1000
$wsNew.Range("E25").Value = "[Two equal peaks at 0.27 and 0.42]"
$wsNew.Range("F25").Value = $wsNew.Range("F24").Value()
$wsNew.Range("G25").Value = 200
$wsNew.Range("H25").Value = 1
$wsNew.Range("I25").Value = 0
$wsNew.Range("J25").Value = 0
$wsNew.Range("K25").Value = 0
$wsNew.Range("L25").Value = 30

# Move the selection on "1D NEW" down past the freshly-added row.
$wsNew.Activate()
$wsNew.Range("E26").Select()

# Leave "1D OLD" as the active tab (matches the saved workbook view state).
$wsOld.Activate()
